# Updates "Price" (D) and "Volume(1h)" (E) columns on the cryptos sheet
# to match the latest scrape, per the GitHub Actions commit.
#
# D-column values that are plain decimal numbers ("123.45") are written
# with a leading apostrophe so Excel keeps storing them as literal text
# (same as the source data / the original cells) instead of silently
# re-typing them as numbers and dropping significant trailing zeros
# (e.g. "0.0500" -> 0.05). Values that already are not valid numbers
# (thousands-dotted prices, the PEPE subscript price, and every
# Volume(1h) percentage string) are assigned as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "56.667.43"
$ws.Range("E2").Value = "  -0.05%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.338.65"
$ws.Range("E3").Value = "  -0.38%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5: BNB
$ws.Range("D5").Value = "'514.51"
$ws.Range("E5").Value = "  -0.36%  "

# Row 6: Solana
$ws.Range("D6").Value = "'133.96"
$ws.Range("E6").Value = "  +0.36%  "

# Row 7: USDC
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.12%  "

# Row 8: XRP
$ws.Range("E8").Value = "  -0.23%  "

# Row 9: Dogecoin
$ws.Range("E9").Value = "  -1.62%  "

# Row 11: Toncoin
$ws.Range("D11").Value = "'5.32"
$ws.Range("E11").Value = "  +1.33%  "

# Row 12: Cardano
$ws.Range("E12").Value = "  -0.17%  "

# Row 13: Avalanche
$ws.Range("D13").Value = "'23.82"
$ws.Range("E13").Value = "  +0.74%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.756.54"
$ws.Range("E14").Value = "  +0.68%  "

# Row 15: WrappedBTC
$ws.Range("D15").Value = "56.645.05"
$ws.Range("E15").Value = "  -0.10%  "

# Row 16: ShibaInu
$ws.Range("E16").Value = "  -0.19%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "2.337.79"
$ws.Range("E17").Value = "  -0.47%  "

# Row 18: Chainlink
$ws.Range("D18").Value = "'10.42"

# Row 19: BitcoinCash
$ws.Range("D19").Value = "'326.40"
$ws.Range("E19").Value = "  +2.27%  "

# Row 20: Polkadot
$ws.Range("D20").Value = "'4.19"
$ws.Range("E20").Value = "  -1.36%  "

# Row 21: Uniswap
$ws.Range("E21").Value = "  +0.82%  "

# Row 22: Dai
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.21%  "

# Row 23: Litecoin
$ws.Range("D23").Value = "'61.14"
$ws.Range("E23").Value = "  +0.83%  "

# Row 24: InternetComputer(DFINITY)
$ws.Range("D24").Value = "'8.72"
$ws.Range("E24").Value = "  +12.98%  "

# Row 25: Kaspa
$ws.Range("E25").Value = "  +3.92%  "

# Row 26: Binance-PegBSC-USD
$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = "  -0.10%  "

# Row 27: Fetch.AI
$ws.Range("E27").Value = "  +7.27%  "

# Row 28: Monero
$ws.Range("D28").Value = "'168.26"
$ws.Range("E28").Value = "  -1.55%  "

# Row 29: PEPE
$ws.Range("D29").Value = "0.0₃0728"
$ws.Range("E29").Value = "  -1.04%  "

# Row 30: PancakeSwap
$ws.Range("D30").Value = "'1.68"
$ws.Range("E30").Value = "  +0.56%  "

# Row 31: Aptos
$ws.Range("D31").Value = "'6.14"
$ws.Range("E31").Value = "  -1.32%  "

# Row 32: EthereumClassic
$ws.Range("D32").Value = "'18.39"
$ws.Range("E32").Value = "  +0.82%  "

# Row 33: USDe
$ws.Range("E33").Value = "  +0.02%  "

# Row 34: FirstDigitalUSD
$ws.Range("E34").Value = "  -0.19%  "

# Row 35: ImmutableX
$ws.Range("E35").Value = "  +2.79%  "

# Row 36: NEARProtocol
$ws.Range("E36").Value = "  +0.62%  "

# Row 37: SuiNetwork
$ws.Range("D37").Value = "'0.889"
$ws.Range("E37").Value = "  -6.16%  "

# Row 38: Stacks
$ws.Range("E38").Value = "  +2.78%  "

# Row 39: OKB
$ws.Range("D39").Value = "'38.64"
$ws.Range("E39").Value = "  +3.36%  "

# Row 40: Aave
$ws.Range("D40").Value = "'150.37"
$ws.Range("E40").Value = "  +9.08%  "

# Row 41: PolygonEcosystemToken
$ws.Range("D41").Value = "'0.376"
$ws.Range("E41").Value = "  -0.63%  "

# Row 42: Filecoin
$ws.Range("E42").Value = "  +1.35%  "

# Row 43: Bittensor
$ws.Range("D43").Value = "'281.34"
$ws.Range("E43").Value = "  +2.17%  "

# Row 44: RenderToken
$ws.Range("E44").Value = "  +1.79%  "

# Row 45: Stellar
$ws.Range("D45").Value = "'0.0925"
$ws.Range("E45").Value = "  -0.28%  "

# Row 46: Hedera
$ws.Range("D46").Value = "'0.0500"
$ws.Range("E46").Value = "  -0.59%  "

# Row 47: Mantle
$ws.Range("D47").Value = "'0.557"
$ws.Range("E47").Value = "  -0.31%  "

# Row 48: InjectiveProtocol
$ws.Range("D48").Value = "'18.30"
$ws.Range("E48").Value = "  +7.42%  "

# Row 49: VeChain
$ws.Range("D49").Value = "'0.0216"
$ws.Range("E49").Value = "  +0.21%  "

# Row 50: EnergySwap
$ws.Range("D50").Value = "'17.10"
$ws.Range("E50").Value = "  +2.15%  "
